$d = $word.ActiveDocument

# 1. Update the table cell text: "Crear perfil paciente" -> "Crear consulta"
$d.Content.Find.Execute("Crear perfil paciente", $false, $false, $false, $false, $false,
                         $true, 1, $false, "Crear consulta", 2)

# 2. Resize the first table's two grid columns (3041 -> 3587 twips, 5781 -> 5235 twips)
#    1 twip = 1/20 point, so divide by 20 to get points for Word's Width property.
$table = $d.Tables.Item(1)
$table.Columns.Item(1).Width = 3587 / 20.0
$table.Columns.Item(2).Width = 5235 / 20.0

# 3. Replace the "sobre los datos faltantes para crear la consulta." text (built from
#    several runs) with a single simplified sentence.
$d.Content.Find.Execute("sobre los datos faltantes para crear la consulta.", $false, $false, $false, $false, $false,
                         $true, 1, $false, "para que se llene uno de los campos.", 2)
